$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species/observation data between row 3 and row 4
$cols = @("A","B","D","E","F","G","H","Q","R")

foreach ($col in $cols) {
    $r3 = $ws.Range($col + "3")
    $r4 = $ws.Range($col + "4")
    $tmp = $r3.Value2
    $r3.Value = $r4.Value2
    $r4.Value = $tmp
}

# The (empty) "Bestämningsmetod" cell also moves from row 4 to row 3
$ws.Range("AF4").Cut($ws.Range("AF3"))
